$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.261.82"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.578.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.78%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.13%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.22%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.578.65"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.84%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.489"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.35%  "

$ws.Range("E10").Value = "  -0.42%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.85"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.16%  "

$ws.Range("E12").Value = "  +0.32%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.187.98"
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000204"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.99%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.576.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.77%  "

$ws.Range("E17").Value = "  +1.52%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.304.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.20%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.45%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "421.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.607"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.16%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.97%  "

$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000119"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.38%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.38%  "

$ws.Range("E28").Value = "  +2.55%  "

$ws.Range("E29").Value = "  -0.30%  "

$ws.Range("E30").Value = "  +0.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.577.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.81%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.156"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.11%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "24.92"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.57%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.43%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.69"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.97%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.53"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.65"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.33%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "174.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0847"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.14"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.877"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.99%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "45.93"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.28%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.82%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.41%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.23%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.53%  "

$ws.Range("E50").Value = "  -5.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.940"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.30%  "
